$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.456557750701904
$ws.Range("B1").Value = 2.768770217895508
$ws.Range("C1").Value = 2.976203203201294
$ws.Range("D1").Value = 3.088655471801758
$ws.Range("E1").Value = 0.8025214672088623
